# Auto update Excel log
# Appends newly-logged sensor events to the PIR, Humidity, Temperature,
# Proximity and mmWave sheets (2026-02-06 sensor sweep).

function Add-LogRows {
    param(
        [string]$SheetName,
        [int]$StartRow,
        [object[]]$Rows
    )

    $ws = $wb.Worksheets.Item($SheetName)

    for ($i = 0; $i -lt $Rows.Count; $i++) {
        $r = $StartRow + $i
        $data = $Rows[$i]

        # Column A holds dates like "2026-02-06" — force text so Excel does
        # not reinterpret the literal as a date serial number, then restore
        # the plain "Normal" style so the cell keeps the sheet's default
        # (unstyled) look, matching every other logged row.
        $ws.Range("A$r").NumberFormat = "@"
        $ws.Range("A$r").Value = $data[0]
        $ws.Range("A$r").Style = "Normal"

        $ws.Range("B$r").Value = $data[1]
        $ws.Range("C$r").Value = $data[2]
        $ws.Range("D$r").Value = $data[3]

        # Column E sometimes holds percentages like "70.0%" — force text so
        # Excel does not reinterpret the literal as a numeric percentage,
        # then restore the default style the same way as column A.
        $ws.Range("E$r").NumberFormat = "@"
        $ws.Range("E$r").Value = $data[4]
        $ws.Range("E$r").Style = "Normal"

        $ws.Range("F$r").Value = $data[5]
    }
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# PIR sheet: rows 165-177
# ---------------------------------------------------------------------------
$pirRows = @(
    @("2026-02-06","09:53:45","09:00","Bathroom","No Motion","Inactive"),
    @("2026-02-06","09:53:48","09:00","Bathroom","No Motion","Inactive"),
    @("2026-02-06","09:53:51","09:00","Bathroom","Motion Detected","Active"),
    @("2026-02-06","09:53:59","09:00","Bathroom","No Motion","Inactive"),
    @("2026-02-06","09:54:00","09:00","Bathroom","Motion Detected","Active"),
    @("2026-02-06","09:54:07","09:00","Bathroom","No Motion","Inactive"),
    @("2026-02-06","09:54:11","09:00","Bathroom","Motion Detected","Active"),
    @("2026-02-06","09:54:19","09:00","Bathroom","No Motion","Inactive"),
    @("2026-02-06","09:54:24","09:00","Bathroom","No Motion","Inactive"),
    @("2026-02-06","09:54:29","09:00","Bathroom","No Motion","Inactive"),
    @("2026-02-06","09:54:30","09:00","Bathroom","Motion Detected","Active"),
    @("2026-02-06","09:54:37","09:00","Bathroom","No Motion","Inactive"),
    @("2026-02-06","09:54:38","09:00","Bathroom","Motion Detected","Active")
)
Add-LogRows "PIR" 165 $pirRows

# ---------------------------------------------------------------------------
# Humidity sheet: rows 92-98
# ---------------------------------------------------------------------------
$humidityRows = @(
    @("2026-02-06","09:53:47","09:00","Bathroom","70.0%","Active"),
    @("2026-02-06","09:53:57","09:00","Bathroom","69.8%","Active"),
    @("2026-02-06","09:54:02","09:00","Bathroom","70.8%","Active"),
    @("2026-02-06","09:54:12","09:00","Bathroom","70.6%","Active"),
    @("2026-02-06","09:54:22","09:00","Bathroom","70.8%","Active"),
    @("2026-02-06","09:54:32","09:00","Bathroom","70.7%","Active"),
    @("2026-02-06","09:54:42","09:00","Bathroom","70.7%","Active")
)
Add-LogRows "Humidity" 92 $humidityRows

# ---------------------------------------------------------------------------
# Temperature sheet: rows 92-98
# ---------------------------------------------------------------------------
$temperatureRows = @(
    @("2026-02-06","09:53:47","09:00","Bathroom","27.8C","Active"),
    @("2026-02-06","09:53:57","09:00","Bathroom","27.7C","Active"),
    @("2026-02-06","09:54:02","09:00","Bathroom","27.8C","Active"),
    @("2026-02-06","09:54:12","09:00","Bathroom","27.8C","Active"),
    @("2026-02-06","09:54:22","09:00","Bathroom","27.8C","Active"),
    @("2026-02-06","09:54:32","09:00","Bathroom","27.8C","Active"),
    @("2026-02-06","09:54:42","09:00","Bathroom","27.8C","Active")
)
Add-LogRows "Temperature" 92 $temperatureRows

# ---------------------------------------------------------------------------
# Proximity sheet: row 8 (single new row — force nested array with the
# unary comma operator so PowerShell doesn't unroll it)
# ---------------------------------------------------------------------------
$proximityRows = @(
    ,@("2026-02-06","09:54:12","09:00","Living Room Main Door","ENTER","User ENTERED Living Room Main Door")
)
Add-LogRows "Proximity" 8 $proximityRows

# ---------------------------------------------------------------------------
# mmWave sheet: row 2 (single new row — same unary comma trick)
# ---------------------------------------------------------------------------
$mmWaveRows = @(
    ,@("2026-02-06","09:54:13","09:00","Living Room Main Door","Image Captured","Active")
)
Add-LogRows "mmWave" 2 $mmWaveRows
